# key_buildingBlock_pairs.xlsx -- "added excel for prob 27 in p&s part"
#
# Sheet1 holds key/building-block triples in columns A (key code),
# B (building block description) and C (supporting expression), one
# problem-step per row. This edit:
#   1. Extends the "$f(p) \times g(p)$;" note in C5 with an extra term.
#   2. Rewords the probability-complement building block in B132.
#   3. Appends three new rows (133-135) describing a new problem
#      (probability 27 in the p&s part): confidence interval for a
#      population mean, solving for the population standard
#      deviation, and a trailing key-only row.
#   4. Moves the viewport/selection to the top of the sheet (C7).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Extend the existing supporting-expression note in C5.
$ws.Range("C5").Value = "`$f(p) \times g(p)`$; `$b-a`$;"

# 2. Reword the building block describing the complement-probability step.
$ws.Range("B132").Value = "사건이 일어날 확률은 그 여사건이 일어날 확률을 계산하고 나서 구합니다."

# 3. Append the new problem's rows.
$ws.Range("A133").Value = "x0013"
$ws.Range("B133").Value = "모집단에서 임의추출하여 구한 표본평균의 값과 그 표본의 크기에 대해 모평균의 신뢰구간을 구합니다."
$ws.Range("C133").Value = "`$a \leq m \leq b`$;"

$ws.Range("A134").Value = "x0014"
$ws.Range("B134").Value = "주어진 조건으로 만들어진 연립방정식을 풀어서 모표준편차를 구합니다."
$ws.Range("C134").Value = "`$0.67 \times \dfrac{\sigma}{10}=1.34`$;"

$ws.Range("A135").Value = "x0015"

# 4. Update the active viewport/selection.
$ws.Activate()
$ws.Range("C7").Select()
